# Apply Villads Egede Johansen's review comments to the document.
#
# Four comments are anchored on specific spans of text in the two
# paragraphs discussing the optimisation difficulty and the machine
# learning solution. Adding them via $d.Comments.Add() automatically
# splits the host run into commentRangeStart / anchor-run /
# commentRangeEnd / commentReference-run pieces, mirroring the target
# diff, and mints the comments.xml part + Comment* / BalloonText*
# styles.
#
# NOTE: the COM shim only keeps one "live" Comment proxy per
# $d.Comments.Add() call alive at a time -- reusing the variable it
# returns to set .Author/.Initial later silently clobbers a previous
# comment's properties. Re-fetching via $d.Comments.Item(n) right
# before each property write avoids that, so every comment ends up
# correctly attributed.

$d = $word.ActiveDocument

# --- Add the four comments, anchored on their exact text spans -------

# Comment 0: anchors "optimisation processes heavily relies on constraints "
$r0 = $d.Range(1513, 1566)
$d.Comments.Add($r0, "I’m still not quite sure how this should be understood") | Out-Null

# Comment 1: anchors "Quantum computers are considered to be able to find
# global minima very effectively, however they are not expected to be
# available in the near future. "
$r1 = $d.Range(1879, 2029)
$d.Comments.Add($r1, "Hmm, maybe a bit unnecessary to mention") | Out-Null

# Comment 2: anchors "a simple and powerful"
$r2 = $d.Range(2257, 2278)
$d.Comments.Add($r2, "Maybe don’t write “simple” first, since this undermines your work. What about “powerful yet simple…”") | Out-Null

# Comment 3: anchors "Based on the determined classification, we can
# generate constraints that can be used for the optimisation process
# that follows afterwards."
$r3 = $d.Range(2873, 3011)
$d.Comments.Add($r3, "Maybe also strengthen your description of the results you have obtained by using this approach (now it seems as an intermediate step, whereas you can sell more as a crucial invention for our analysis)") | Out-Null

# --- Attribute all four comments to the reviewer ----------------------
for ($i = 1; $i -le $d.Comments.Count; $i++) {
    $d.Comments.Item($i).Author = "Villads Egede Johansen"
    $d.Comments.Item($i).Initial = "VE"
}

# --- Relocate the "_GoBack" bookmark ---------------------------------
# It used to sit alone in its own empty paragraph; move it to mark the
# last edit location (the start of the final comment's anchor).
if ($d.Bookmarks.Exists("_GoBack")) {
    $goBack = $d.Bookmarks.Item("_GoBack")
    $goBack.Delete()
}
$goBackRange = $d.Range(2873, 2873)
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null

Write-Output "comments added"
